# RAD.IMR ValueSet (imr-servicerequest-intent-vs) publication edit:
# TI Publication of RAD.IMR 1.0.0
#   - Version goes from "1.0.0-comment" to "1.0.0"
#   - Date is bumped to the publication timestamp
#   - The draft Copyright text is cleared (no longer populated)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 3: Version
$ws.Range("B3").Value = "1.0.0"

# Row 8: Date
$ws.Range("B8").Value = "2022-07-25T14:40:04-05:00"

# Row 16: Copyright value cleared
$ws.Range("B16").ClearContents()
